$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chillers: Amount 2 -> 1
$ws.Range("D2").Value = 1

# Pumps: Width 2 -> 0, Length 2 -> 0
$ws.Range("M3").Value = 0
$ws.Range("M4").Value = 0

# Update the active selection to G3
$ws.Range("G3").Select()
